$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of 1-based row index -> new cell(1) text, matching the XML diff.
$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "44"
    6  = "0.00037"
    7  = "0.00016"
    8  = "0.00004"
    9  = "0.00021"
    10 = "0.00026"
    11 = "0.00029"
    12 = "0.00704"
    44 = "99.94"
    45 = "0.01"
    46 = "12"
}

foreach ($rowIndex in $updates.Keys) {
    $cell = $t.Cell($rowIndex, 1)
    $cell.Range.Text = $updates[$rowIndex]
}
